# CodingQuestionsHint.xlsx - "Add files via upload" edit
#
# Adds two new rows of coding-question tracker data (rows 65 & 66) to
# Sheet1, reusing the existing column layout / styles used by the other
# rows in the table, and nudges the sheet view / row heights to line up
# with the refreshed content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New row 65: "Remove Outermost Parentheses" (Leetcode)
# ---------------------------------------------------------------------
$ws.Range("A65").Value = 63

# B65 - date (reuse the date-formatted style already used by column B)
$ws.Range("B64").Copy()
$ws.Range("B65").PasteSpecial(-4122)
$ws.Range("B65").Value = 45733

# C65 - problem title, partially bold ("Leetcode" suffix), no-fill style
$ws.Range("C65").Style = "Normal"
$ws.Range("C65").Value = "Remove Outermost Parentheses. Leetcode"
$leet65 = $ws.Range("C65").Characters(31, 8)
$leet65.Font.Bold = $true
$leet65.Font.ColorIndex = -4105
$ws.Range("C65").Interior.ColorIndex = -4142

# D65 - input (reuse the wrap-text style used elsewhere in column D)
$ws.Range("C64").Copy()
$ws.Range("D65").PasteSpecial(-4122)
$ws.Range("D65").Value = 's = "(()())(())"'

# E65 - output (reuse the Consolas "code" style used elsewhere in column E)
$ws.Range("E39").Copy()
$ws.Range("E65").PasteSpecial(-4122)
$ws.Range("E65").Value = '"()()()"'

# F65 - solution note (reuse existing shared string "O(n) - Easy")
$ws.Range("F64").Copy()
$ws.Range("F65").PasteSpecial(-4122)
$ws.Range("F65").Value = "O(n) - Easy"

# G65 - who (reuse existing shared string "took help")
$ws.Range("G64").Copy()
$ws.Range("G65").PasteSpecial(-4122)
$ws.Range("G65").Value = "took help"

# ---------------------------------------------------------------------
# 2. New row 66: "Longest Common Prefix" (Leetcode)
# ---------------------------------------------------------------------
$ws.Range("A66").Value = 64

# B66 - same date as row 65
$ws.Range("B64").Copy()
$ws.Range("B66").PasteSpecial(-4122)
$ws.Range("B66").Value = 45733

# C66 - problem title (plain text this time), same no-fill style as C65
$ws.Range("C65").Copy()
$ws.Range("C66").PasteSpecial(-4122)
$ws.Range("C66").Value = " Longest Common Prefix"

# D66 - input (Consolas "code" style)
$ws.Range("E39").Copy()
$ws.Range("D66").PasteSpecial(-4122)
$ws.Range("D66").Value = 'strs = ["flower","flow","flight"]'

# E66 - output (Consolas "code" style)
$ws.Range("E39").Copy()
$ws.Range("E66").PasteSpecial(-4122)
$ws.Range("E66").Value = '"fl"'

# F66 - solution note (reuse existing shared string "O(n) - Easy")
$ws.Range("F64").Copy()
$ws.Range("F66").PasteSpecial(-4122)
$ws.Range("F66").Value = "O(n) - Easy"

# G66 - who (reuse existing shared string "solved and submitted")
$ws.Range("C35").Copy()
$ws.Range("G66").PasteSpecial(-4122)
$ws.Range("G66").Value = "solved and submitted"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Row-height refresh (content/line-height recalculation across the
#    table caused by the new rows).
# ---------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 60
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 60
$ws.Rows.Item(6).RowHeight = 60
$ws.Rows.Item(7).RowHeight = 60
$ws.Rows.Item(8).RowHeight = 60
$ws.Rows.Item(9).RowHeight = 60
$ws.Rows.Item(10).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(12).RowHeight = 60
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 60
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 30
$ws.Rows.Item(22).RowHeight = 90
$ws.Rows.Item(23).RowHeight = 210
$ws.Rows.Item(25).RowHeight = 60
$ws.Rows.Item(26).RowHeight = 60
$ws.Rows.Item(27).RowHeight = 60
$ws.Rows.Item(29).RowHeight = 120
$ws.Rows.Item(30).RowHeight = 60
$ws.Rows.Item(31).RowHeight = 75
$ws.Rows.Item(33).RowHeight = 45
$ws.Rows.Item(35).RowHeight = 60
$ws.Rows.Item(36).RowHeight = 60
$ws.Rows.Item(38).RowHeight = 60
$ws.Rows.Item(39).RowHeight = 60
$ws.Rows.Item(40).RowHeight = 60
$ws.Rows.Item(41).RowHeight = 30
$ws.Rows.Item(42).RowHeight = 30
$ws.Rows.Item(43).RowHeight = 30
$ws.Rows.Item(44).RowHeight = 60
$ws.Rows.Item(45).RowHeight = 45
$ws.Rows.Item(46).RowHeight = 75
$ws.Rows.Item(47).RowHeight = 60
$ws.Rows.Item(48).RowHeight = 60
$ws.Rows.Item(49).RowHeight = 60
$ws.Rows.Item(50).RowHeight = 30
$ws.Rows.Item(51).RowHeight = 60
$ws.Rows.Item(52).RowHeight = 60
$ws.Rows.Item(53).RowHeight = 60
$ws.Rows.Item(54).RowHeight = 60
$ws.Rows.Item(55).RowHeight = 60
$ws.Rows.Item(56).RowHeight = 60
$ws.Rows.Item(57).RowHeight = 60
$ws.Rows.Item(58).RowHeight = 60
$ws.Rows.Item(59).RowHeight = 60
$ws.Rows.Item(60).RowHeight = 60
$ws.Rows.Item(61).RowHeight = 30
$ws.Rows.Item(62).RowHeight = 60
$ws.Rows.Item(63).RowHeight = 60
$ws.Rows.Item(64).RowHeight = 30
$ws.Rows.Item(65).RowHeight = 30
$ws.Rows.Item(66).RowHeight = 60

# ---------------------------------------------------------------------
# 4. Sheet view: scroll/select near the new rows (where the author's
#    cursor ended up after the edit).
# ---------------------------------------------------------------------
$ws.Range("E66").Select()
